$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: J4 text update
$ws.Range("J4").Value = "Avaliable by the end of the year"

# Row 5: shorten title in A5
$ws.Range("A5").Value = "Many labs 3"

# Row 6: shift old full citation (A6) into B6, set short title into A6, add E6/F6
$ws.Range("B6").Value = $ws.Range("A6").Value()
$ws.Range("A6").Value = "Evaluating the replicability of social science experiments in Nature and Science between 2010 and 2015"
$ws.Range("E6").Value = " socialScienceExperimentsInNatureAndScience.csv"
$ws.Range("F6").Value = "Check tables at e.g., https://osf.io/bh9xs/ to verify what each column is"

# Row 7: shift old full citation (A7) into B7, set short title into A7, update C7, clear D7, add E7
$ws.Range("B7").Value = $ws.Range("A7").Value()
$ws.Range("A7").Value = "Evaluating replicability of laboratory experiments in economics"
$ws.Range("C7").Value = "https://osf.io/bzm54/"
$ws.Range("D7").ClearContents()
$ws.Range("E7").Value = "socialScienceExperimentsInNatureAndScience.csv"

# Remove row 10 (the SUM formula row) entirely
$ws.Rows.Item(10).Delete()

# Update selection to match final state
$ws.Range("J4").Select()
